$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 292.7353
$ws.Range("I6").Value = 294.93103
$ws.Range("J6").Value = 280
$ws.Range("K6").Value = 884.7930900000001
$ws.Range("L6").Value = 840
$ws.Range("M6").Value = -772.7930900000001
$ws.Range("N6").Value = -1064
$ws.Range("H8").Value = 99.5
$ws.Range("I8").Value = 99.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 298.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -159.5
$ws.Range("N8").ClearContents()
$ws.Range("H17").Value = 2833872
$ws.Range("J17").Value = 3069986.2
$ws.Range("L17").Value = 9209958.600000001
$ws.Range("N17").Value = -9210294.600000001
$ws.Range("H64").Value = 1879701.2
$ws.Range("I64").Value = 3348081
$ws.Range("J64").Value = 3438.3333
$ws.Range("K64").Value = 3348081
$ws.Range("L64").Value = 3438.3333
$ws.Range("M64").Value = -3347833
$ws.Range("N64").Value = -3934.3333
$ws.Range("H67").Value = 1879701.2
$ws.Range("I67").Value = 3348081
$ws.Range("J67").Value = 3438.3333
$ws.Range("K67").Value = 3348081
$ws.Range("L67").Value = 3438.3333
$ws.Range("M67").Value = -3347223
$ws.Range("N67").Value = -5154.3333
$ws.Range("H132").Value = 7582483
$ws.Range("I132").Value = 6256.7036
$ws.Range("J132").Value = 41675500
$ws.Range("K132").Value = 18770.1108
$ws.Range("L132").Value = 125026500
$ws.Range("M132").Value = -16240.1108
$ws.Range("N132").Value = -125031560
$ws.Range("H138").Value = 4018410.2
$ws.Range("I138").Value = 7937928.5
$ws.Range("J138").Value = 3294.0732
$ws.Range("K138").Value = 23813785.5
$ws.Range("L138").Value = 9882.2196
$ws.Range("M138").Value = -23808645.5
$ws.Range("N138").Value = -20162.2196

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10663.667
$ws.Range("I3").Value = 1005
$ws.Range("J3").Value = 15493
$ws.Range("K3").Value = 1005
$ws.Range("L3").Value = 15493
$ws.Range("M3").Value = -890
$ws.Range("N3").Value = -15723
$ws.Range("H32").Value = 6897
$ws.Range("I32").Value = 5570.6294
$ws.Range("J32").Value = 17628.545
$ws.Range("K32").Value = 5570.6294
$ws.Range("L32").Value = 17628.545
$ws.Range("M32").Value = -5283.6294
$ws.Range("N32").Value = -18202.545
$ws.Range("H35").Value = 1714.5
$ws.Range("I35").Value = 939
$ws.Range("J35").Value = 4041
$ws.Range("K35").Value = 939
$ws.Range("L35").Value = 4041
$ws.Range("M35").Value = -533
$ws.Range("N35").Value = -4853
$ws.Range("H63").Value = 83335200
$ws.Range("I63").Value = 100001660
$ws.Range("J63").Value = 2900
$ws.Range("K63").Value = 100001660
$ws.Range("L63").Value = 2900
$ws.Range("M63").Value = -100000974
$ws.Range("N63").Value = -4272
$ws.Range("H66").Value = 83335200
$ws.Range("I66").Value = 100001660
$ws.Range("J66").Value = 2900
$ws.Range("K66").Value = 500008300
$ws.Range("L66").Value = 14500
$ws.Range("M66").Value = -500004868
$ws.Range("N66").Value = -21364
$ws.Range("H74").Value = 16668843
$ws.Range("I74").Value = 23810834
$ws.Range("J74").Value = 4198.1113
$ws.Range("K74").Value = 23810834
$ws.Range("L74").Value = 4198.1113
$ws.Range("M74").Value = -23809960
$ws.Range("N74").Value = -5946.1113
$ws.Range("H77").Value = 16668843
$ws.Range("I77").Value = 23810834
$ws.Range("J77").Value = 4198.1113
$ws.Range("K77").Value = 119054170
$ws.Range("L77").Value = 20990.5565
$ws.Range("M77").Value = -119049802
$ws.Range("N77").Value = -29726.5565

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1492
$ws.Range("I11").Value = 1492
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1492
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1352
$ws.Range("N11").ClearContents()
$ws.Range("H81").Value = 21585
$ws.Range("J81").Value = 21585
$ws.Range("L81").Value = 21585
$ws.Range("N81").Value = -23707
$ws.Range("H84").Value = 21585
$ws.Range("J84").Value = 21585
$ws.Range("L84").Value = 64755
$ws.Range("N84").Value = -75363

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 22730880
$ws.Range("I132").Value = 33336424
$ws.Range("J132").Value = 4718
$ws.Range("K132").Value = 100009272
$ws.Range("L132").Value = 14154
$ws.Range("M132").Value = -100006742
$ws.Range("N132").Value = -19214

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6666853.5
$ws.Range("I4").Value = 6666853.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 20000560.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -20000448.5
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 300
$ws.Range("J7").Value = 300
$ws.Range("L7").Value = 900
$ws.Range("N7").Value = -1124
$ws.Range("H68").Value = 975.64703
$ws.Range("I68").Value = 1248.4
$ws.Range("J68").Value = 862
$ws.Range("K68").Value = 3745.2
$ws.Range("L68").Value = 2586
$ws.Range("M68").Value = -2934.2
$ws.Range("N68").Value = -4208
$ws.Range("H71").Value = 975.64703
$ws.Range("I71").Value = 1248.4
$ws.Range("J71").Value = 862
$ws.Range("K71").Value = 11235.6
$ws.Range("L71").Value = 7758
$ws.Range("M71").Value = -7179.6
$ws.Range("N71").Value = -15870
$ws.Range("H120").Value = 12981
$ws.Range("I120").Value = 2943.3333
$ws.Range("J120").Value = 17999.834
$ws.Range("K120").Value = 8829.999899999999
$ws.Range("L120").Value = 53999.50199999999
$ws.Range("M120").Value = -3991.999899999999
$ws.Range("N120").Value = -63675.50199999999
$ws.Range("H121").Value = 457
$ws.Range("I121").Value = 430.55554
$ws.Range("J121").Value = 933
$ws.Range("K121").Value = 1291.66662
$ws.Range("L121").Value = 2799
$ws.Range("M121").Value = 18.33338000000003
$ws.Range("N121").Value = -5419

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
$ws.Range("H133").Value = 66171.92999999999
$ws.Range("J133").Value = 66171.92999999999
$ws.Range("L133").Value = 66171.92999999999
$ws.Range("N133").Value = -76291.92999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3963.4595
$ws.Range("I40").Value = 4962
$ws.Range("J40").Value = 3202.6667
$ws.Range("K40").Value = 4962
$ws.Range("L40").Value = 3202.6667
$ws.Range("M40").Value = -4826
$ws.Range("N40").Value = -3474.6667
$ws.Range("H61").Value = 2071.25
$ws.Range("I61").Value = 1945
$ws.Range("J61").Value = 2450
$ws.Range("K61").Value = 1945
$ws.Range("L61").Value = 2450
$ws.Range("M61").Value = -1743
$ws.Range("N61").Value = -2854
$ws.Range("H68").Value = 2914.4443
$ws.Range("J68").Value = 3200
$ws.Range("L68").Value = 3200
$ws.Range("N68").Value = -4698
$ws.Range("H71").Value = 2914.4443
$ws.Range("J71").Value = 3200
$ws.Range("L71").Value = 16000
$ws.Range("N71").Value = -23488
$ws.Range("H82").Value = 3057.1428
$ws.Range("I82").Value = 2966.6667
$ws.Range("J82").Value = 3125
$ws.Range("K82").Value = 2966.6667
$ws.Range("L82").Value = 3125
$ws.Range("M82").Value = -2605.6667
$ws.Range("N82").Value = -3847
$ws.Range("H85").Value = 3057.1428
$ws.Range("I85").Value = 2966.6667
$ws.Range("J85").Value = 3125
$ws.Range("K85").Value = 2966.6667
$ws.Range("L85").Value = 3125
$ws.Range("M85").Value = -1718.6667
$ws.Range("N85").Value = -5621
$ws.Range("H113").Value = 2071.25
$ws.Range("I113").Value = 1945
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 1945
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = 225
$ws.Range("N113").Value = -6790
$ws.Range("H122").Value = 6365.1
$ws.Range("I122").Value = 8250.5
$ws.Range("J122").Value = 5108.1665
$ws.Range("K122").Value = 24751.5
$ws.Range("L122").Value = 15324.4995
$ws.Range("M122").Value = -22301.5
$ws.Range("N122").Value = -20224.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8900
$ws.Range("I4").Value = 8000
$ws.Range("J4").Value = 9800
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 9800
$ws.Range("M4").Value = -7887
$ws.Range("N4").Value = -10026
$ws.Range("H62").Value = 6128
$ws.Range("I62").Value = 5633.3335
$ws.Range("J62").Value = 6870
$ws.Range("K62").Value = 5633.3335
$ws.Range("L62").Value = 6870
$ws.Range("M62").Value = -5009.3335
$ws.Range("N62").Value = -8118
$ws.Range("H65").Value = 6128
$ws.Range("I65").Value = 5633.3335
$ws.Range("J65").Value = 6870
$ws.Range("K65").Value = 28166.6675
$ws.Range("L65").Value = 34350
$ws.Range("M65").Value = -25046.6675
$ws.Range("N65").Value = -40590
$ws.Range("H86").Value = 1960.7142
$ws.Range("J86").Value = 1960.7142
$ws.Range("L86").Value = 1960.7142
$ws.Range("N86").Value = -4206.7142
$ws.Range("H89").Value = 1960.7142
$ws.Range("J89").Value = 1960.7142
$ws.Range("L89").Value = 9803.571
$ws.Range("N89").Value = -21035.571
$ws.Range("H113").Value = 414.91666
$ws.Range("I113").Value = 131.33333
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 393.99999
$ws.Range("L113").Value = 7200
$ws.Range("M113").Value = 1776.00001
$ws.Range("N113").Value = -11540
$ws.Range("H135").Value = 29943
$ws.Range("J135").Value = 29943
$ws.Range("L135").Value = 29943
$ws.Range("N135").Value = -40083
$ws.Range("H136").Value = 878.7954999999999
$ws.Range("I136").Value = 946.0833
$ws.Range("J136").Value = 576
$ws.Range("K136").Value = 2838.2499
$ws.Range("L136").Value = 1728
$ws.Range("M136").Value = -288.2498999999998
$ws.Range("N136").Value = -6828
